$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Column AS holds the "subgenus" header (row 1) / "${subgenus}" template
# value (row 2), which review determined should not be part of the
# Materials mapping sheet. Remove the entire column, shifting everything
# to its right one column to the left.
$ws.Columns("AS:AS").Delete()
